$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make A2 bold (reuses existing bold Times New Roman style already in the workbook)
$ws.Range("A2").Font.Bold = $true

# Update the raw parameter values (B2:D4) with the new "deaths at each step" figures
# and clear their prior numeric style so they fall back to the default style,
# matching the target workbook.
$values = @{
    "B2" = 9.761902954
    "C2" = 2.2706403989999999
    "D2" = 0.29226176269999998
    "B3" = 6.6876464320000002
    "C3" = 10.333415990000001
    "D3" = 0.88528063099999998
    "B4" = 0.91568368420000001
    "C4" = 1.01411684
    "D4" = 1.2494830159999999
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).ClearFormats()
    $ws.Range($addr).Value = $values[$addr]
}

# Update the current selection to match the new working range
$ws.Range("A15:E20").Select()
